# Apply the "Updated symbol list" refresh: prices in column D, the
# "Worstin24h"/"Bestin24h" badge suffix that jumps between rows in column E,
# and the snapshot hour in column G (11 -> 12) for every data row (2-51).
# Numeric-looking values are written with a leading apostrophe so Excel
# keeps storing them as text, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.62"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'23.89"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.259"
$ws.Range("G4").Value = "'12"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.462"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'3.331"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.8099"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8757"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.07266"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.03087"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.03055"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.09324"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'3.847"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.001539"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'0.04719"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'0.0006037"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'0.006144"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.001263"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.004599"
$ws.Range("G21").Value = "'12"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'3.562"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'2.181"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.3209"
$ws.Range("G25").Value = "'12"
$ws.Range("G26").Value = "'12"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03781"
$ws.Range("G40").Value = "'12"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.1051"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.002329"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.007302"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005476"
$ws.Range("G45").Value = "'12"
$ws.Range("G46").Value = "'12"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.01856"
$ws.Range("G48").Value = "'12"
$ws.Range("G49").Value = "'12"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
